# Automated map update (2025-09-02 08:02:14)
# A new case (Caso 4768) was logged on 1/30/2025, which sorts before the
# existing row 25 (Caso 6045, 2/7/2025). Insert a fresh row at position 25
# so every following record shifts down by one, then populate the new row
# with the new case's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Insert()

# Columns A, B, D and E hold numeric/date-looking text (Caso, F. De
# Reclamo, Comuna, OT) that must stay plain text, matching every other
# row in the sheet - force the Text number format first so Excel doesn't
# silently coerce them into a number/date.
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "4768"

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "1/30/2025"

$ws.Range("C25").Value = "VALLESE, FELIPE 684"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "802988221"

$ws.Range("F25").Value = "PEBCOM"
$ws.Range("G25").Value = "Pendiente"
$ws.Range("H25").Value = "Picada info para cierre tambien entro como caso 6909"
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = "Cambio"
$ws.Range("K25").Value = "Sin equipos"
$ws.Range("L25").Value = "Pasante"
$ws.Range("M25").Value = -58.443039
$ws.Range("N25").Value = -34.612262
$ws.Range("O25").Value = "Almagro"
$ws.Range("P25").Value = "Capital Sur"
